$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Python")
Write-Host $ws.Name
